$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 63: new LeetCode entry - 876: Middle of the Linked List
$ws.Range("A63").Value = 876
$ws.Range("B63").Value = "Middle of the Linked List"
$ws.Range("C63").Value = "#linked-list #two-pointers #重点 "
$ws.Range("D63").Value = "easy"
$ws.Range("E63").Value = 0
$ws.Range("F63").Value = 3
$ws.Range("G63").Value = 1
$ws.Range("H63").Value = [DateTime]"2025-07-08"
$ws.Range("I63").Value = [DateTime]"2025-07-08"

# Row 64: extra tags line associated with the entry above (Tags column only)
$ws.Range("C64").Value = "#array #sorting #核心 "

$ws.Rows("63").RowHeight = 34
$ws.Rows("64").RowHeight = 34

$ws.Range("B64").Select()
